$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 50 ("「腹空き兎」...") entirely; subsequent rows shift up by one.
$ws.Rows.Item(50).Delete()
